$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Replace the "NA" dispensador values with "Tork® Jabones" for the
# "Jabones y Gel" rows (rows 39-52, column D).
for ($r = 39; $r -le 52; $r++) {
    $ws.Cells.Item($r, 4).Value = "Tork® Jabones"
}

# Extend the AutoFilter range to include the new last row (A1:F52).
$ws.AutoFilterMode = $false
$ws.Range("A1:F52").AutoFilter()

# Keep the workbook-level _FilterDatabase defined name in sync with the
# AutoFilter range.
foreach ($n in $wb.Names) {
    if ($n.Name -eq "Hoja1!_FilterDatabase") {
        $n.RefersTo = "=Hoja1!`$A`$1:`$F`$52"
    }
}

# Update the active selection on the sheet.
$ws.Range("C43").Select()
